$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 81, shifting existing rows 81-105 down to 82-106.
$ws.Rows(81).Insert()

# Populate the newly inserted row 81 with the new weekly record.
$ws.Cells.Item(81, 1).Value = 5
$ws.Cells.Item(81, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(81, 3).Value = "Maule"
$ws.Cells.Item(81, 4).Value = 44508
$ws.Cells.Item(81, 5).Value = 7
$ws.Cells.Item(81, 6).Value = 100112031
$ws.Cells.Item(81, 7).Value = "Poroto verde"
$ws.Cells.Item(81, 8).Value = "Sin especificar"
$ws.Cells.Item(81, 9).Value = "Primera"
$ws.Cells.Item(81, 10).Value = 150
$ws.Cells.Item(81, 11).Value = 40000
$ws.Cells.Item(81, 12).Value = 40000
$ws.Cells.Item(81, 13).Value = 40000
$ws.Cells.Item(81, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(81, 15).Value = "Región del Maule"
$ws.Cells.Item(81, 16).Value = 1600
$ws.Cells.Item(81, 17).Value = 25
$ws.Cells.Item(81, 18).Value = "Hortaliza"
